$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sits inside the TOC entry
#    "Scope and Pu|rpose" and merge the two split runs back into one.
#    Deleting the bookmark first (while the two runs are still split)
#    then doing a targeted Find/Replace over the phrase collapses the
#    surrounding runs (identical rPr) back into a single run, and
#    causes every other bookmark's w:id to be renumbered sequentially
#    from 0 in document order -- exactly what the target XML expects.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng = $d.Content
$rng.Find.Execute("Scope and Purpose", $true, $false, $false, $false, $false, $true, 1, $false, "Scope and Purpose", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Rewrite the two bullet paragraphs under the "Update Error Types"
#    heading.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {

    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "User must have loaded the error type into the Error Details form.") {

        $r = $p.Range
        $r.SetRange($r.Start, $r.End - 1)
        $r.Text = ""
        $ins = $d.Range($r.Start, $r.Start)
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">User must have </w:t></w:r><w:r><w:t>located</w:t></w:r><w:r><w:t xml:space="preserve"> the error type </w:t></w:r><w:r><w:t>via Search</w:t></w:r><w:r><w:t xml:space="preserve"> Error </w:t></w:r><w:r><w:t>Types</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $ins.InsertXML($xml)
    }
}

foreach ($p in $d.Paragraphs) {

    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "User must know what the field is to be updated to.") {

        # A fresh "_GoBack" bookmark is (re)created around "the " -- Word
        # renumbers every bookmark's w:id by document position on save, so
        # the literal id used here ("900") doesn't matter.
        $r = $p.Range
        $r.SetRange($r.Start, $r.End - 1)
        $r.Text = ""
        $ins = $d.Range($r.Start, $r.Start)
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">User must know what </w:t></w:r><w:bookmarkStart w:id="900" w:name="_GoBack"/><w:r><w:t xml:space="preserve">the </w:t></w:r><w:bookmarkEnd w:id="900"/><w:r><w:t>field is to be updated to</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $ins.InsertXML($xml)
    }
}
